# Powerpoint writer: consolidate text runs when possible.
#
# The shapes below currently hold their text split across many
# single-word <a:r> runs (one run per word/space). Re-assigning the
# same final text in one shot would be a no-op for the COM host's
# diffing (it only rewrites runs that actually changed), so each
# TextRange is first set to a short placeholder value to invalidate
# the existing run split, then set to the desired consolidated text.
# This collapses the paragraph down to a single run that reuses the
# first run's (empty) <a:rPr/>, matching the "slimmed down" output.

function Consolidate-Text($textRange, $finalText) {
    $textRange.Text = "x"
    $textRange.Text = $finalText
}

$p = $ppt.ActivePresentation

# Slide 1: "Section Header (with background image)" title.
$s1 = $p.Slides.Item(1)
Consolidate-Text $s1.Shapes.Item(1).TextFrame.TextRange "Section Header (with background image)"

# Slide 2: "Slide 1" title.
$s2 = $p.Slides.Item(2)
Consolidate-Text $s2.Shapes.Item(1).TextFrame.TextRange "Slide 1"

# Slide 3: "Slide 2" title.
$s3 = $p.Slides.Item(3)
Consolidate-Text $s3.Shapes.Item(1).TextFrame.TextRange "Slide 2"

# Slide 4: "Slide 3" title.
$s4 = $p.Slides.Item(4)
Consolidate-Text $s4.Shapes.Item(1).TextFrame.TextRange "Slide 3"

# Slide 5: "Slide 4" title and the "An image" caption textbox.
$s5 = $p.Slides.Item(5)
Consolidate-Text $s5.Shapes.Item(1).TextFrame.TextRange "Slide 4"
Consolidate-Text $s5.Shapes.Item(4).TextFrame.TextRange "An image"

# Slide 6's notes page: "Blank slides can have background images."
$s6 = $p.Slides.Item(6)
$notesPage = $s6.NotesPage
Consolidate-Text $notesPage.Shapes.Item(2).TextFrame.TextRange "Blank slides can have background images."
